$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("penet")

$ws.Range("A2").Value = "Drilling rig"
$ws.Range("A3").Value = "Hammer"
$ws.Range("A4").Value = "Vibro driver"
$ws.Range("A5").Value = "ROV with suction pump"
$ws.Range("A6").Value = "ROV with jetting"
